$d = $word.ActiveDocument

# 1. Update business name
$d.Content.Find.Execute("Name of the Dog Grooming Business: Puppy Bath Central", $true, $false, $false, $false, $false, $true, 1, $false, "Name of the Dog Grooming Business: The Regal Hound", 2)

# 2. Update Idea 1-5 headings
$d.Content.Find.Execute("Idea 1: ", $true, $false, $false, $false, $false, $true, 1, $false, "Idea 1: About and Credentials", 2)
$d.Content.Find.Execute("Idea 2:", $true, $false, $false, $false, $false, $true, 1, $false, "Idea 2: Services offered", 2)
$d.Content.Find.Execute("Idea 3: ", $true, $false, $false, $false, $false, $true, 1, $false, "Idea 3: Cost for services", 2)
$d.Content.Find.Execute("Idea 4:", $true, $false, $false, $false, $false, $true, 1, $false, "Idea 4: Before & After Gallery", 2)
$d.Content.Find.Execute("Idea 5:", $true, $false, $false, $false, $false, $true, 1, $false, "Idea 5: Appointment booking and policies", 2)
